$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 28
$ws.Cells.Item($row, 1).Value = "Wav2Vec2 Phoneme"
$ws.Cells.Item($row, 2).Value = "Wav2Vec2"
$ws.Cells.Item($row, 3).Value = "Aug_Comb"
$ws.Cells.Item($row, 4).Value = 34.6
$ws.Cells.Item($row, 5).Value = 65.40000000000001
$ws.Cells.Item($row, 6).Value = 0.5766
$ws.Cells.Item($row, 7).Value = 0.6059
$ws.Cells.Item($row, 8).Value = 0.8464
$ws.Cells.Item($row, 9).Value = 0.4904
$ws.Cells.Item($row, 10).Value = 0.8512999999999999
$ws.Cells.Item($row, 11).Value = 0.3444
$ws.Cells.Item($row, 12).Value = 0.49
$ws.Cells.Item($row, 13).Value = " "
